# 2014-lecture1-welcome.pptx — "updated screenshots /ht aswathy"
#
# 1) Slide 24 ("Code of Conduct"): the body placeholder gets a new centered
#    paragraph appended after the "I will post Judi Brown Clark's contact
#    information..." line, and the text frame is switched on to shrink text
#    on overflow (PowerPoint adds <a:normAutofit/> once the extra line no
#    longer fits the placeholder).
# 2) Slide 39 ("Communication!"): the title is retitled to
#    "Process and materials!".

$p = $ppt.ActivePresentation

# --- Slide 24: Code of Conduct -------------------------------------------
$s24 = $p.Slides.Item(24)
$body = $s24.Shapes.Item("Content Placeholder 2")

$tr = $body.TextFrame.TextRange
$para1 = $tr.InsertAfter("`rNote: this is not because of known prior ")
$para2 = $tr.InsertAfter("problems, ICYW.")

# The extra line pushes the placeholder's text past its bounds, so
# PowerPoint shrinks the text to fit (Format Shape > Text Options >
# Shrink text on overflow == ppAutoSizeTextToFitShape).
$body.TextFrame.AutoSize = 2

# --- Slide 39: Communication! ---------------------------------------------
$s39 = $p.Slides.Item(39)
$title = $s39.Shapes.Item("Title 1")
$title.TextFrame.TextRange.Text = "Process and materials!"
